$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) and Volume(1h) (E) columns so numeric-looking
# strings with multiple dots / percent signs / padding spaces are preserved as text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "35.264.12"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.877.51"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "246.31"
$ws.Range("E5").Value = "  -3.07%  "
$ws.Range("D6").Value = "0.680"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").Value = "43.56"
$ws.Range("E8").Value = "  +4.17%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "53.57"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "0.0738"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "0.0978"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "13.50"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "2.150.88"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "0.766"
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "1.895.38"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "35.323.02"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "72.57"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "0.0₃0821"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "243.64"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "12.82"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "4.97"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("E24").Value = "  +7.67%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -5.06%  "
$ws.Range("D27").Value = "165.60"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "8.56"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "18.25"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "2.05"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("B32").Value = "TrustWalletToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D32").Value = "1.68"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").Value = "4.29"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "0.0592"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").Value = "0.841"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").Value = "0.0725"
$ws.Range("E39").Value = "  +10.09%  "
$ws.Range("D40").Value = "17.57"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "96.10"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").Value = "1.304.15"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "0.0804"
$ws.Range("E46").Value = "  +6.70%  "
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").Value = "2.72"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "11.86"
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("E50").Value = "  -5.78%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.056.35"
$ws.Range("E51").Value = "  -1.67%  "
